# Append a new data row (row 63) to each of the 4 sensor-log worksheets.
# Each sheet already has data rows through row 62; this adds one more
# reading (row 63) following the same per-sheet pattern, one hour later
# than the prior row, with an updated G (ID_DEC) value.

$wb = $excel.ActiveWorkbook

$rowsData = @(
    @{
        Sheet = "ROW35-FE-LIFTER"
        A = "2025-03-06 22:42:06"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x d"
        F = 400
        G = "568631262647113770877196"
        H = 400
        I = 13
    },
    @{
        Sheet = "ROW35-MID-LIFTER"
        A = "2025-03-06 22:29:35"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0x e"
        F = 400
        G = "568631262647113770942732"
        H = 400
        I = 14
    },
    @{
        Sheet = "ROW02-FE-LIFTER"
        A = "2025-03-06 22:51:45"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
        D = "0x01,0x90,"
        E = "0xff"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 255
    },
    @{
        Sheet = "ROW02-MID-LIFTER"
        A = "2025-03-06 22:41:15"
        B = "0x01,0x90 "
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x01,0x90,"
        E = "0x 3"
        F = 400
        G = "568631262647113769959692"
        H = 400
        I = 3
    }
)

foreach ($row in $rowsData) {
    $ws = $wb.Worksheets.Item($row.Sheet)
    $newRow = 63

    # Text columns (A-E). Column G is a long digit string that must stay
    # text (it exceeds numeric precision), so force a text quote-prefix on
    # it before assigning the value.
    $ws.Cells.Item($newRow, 1).Value = $row.A
    $ws.Cells.Item($newRow, 2).Value = $row.B
    $ws.Cells.Item($newRow, 3).Value = $row.C
    $ws.Cells.Item($newRow, 4).Value = $row.D
    $ws.Cells.Item($newRow, 5).Value = $row.E

    # Numeric columns.
    $ws.Cells.Item($newRow, 6).Value = $row.F

    # Column G: keep as text (the value is a 24-digit id, not a number).
    $ws.Cells.Item($newRow, 7).Value = "'" + $row.G

    $ws.Cells.Item($newRow, 8).Value = $row.H
    $ws.Cells.Item($newRow, 9).Value = $row.I
}
